# Swap the data values of row 5 and row 6 (columns A:E) on the active sheet.
# (Note: this runtime's Range.Value getter is unreliable, so Value2 is used
# for both reading and writing cell data.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5 = $ws.Range("A5:E5").Value2
$row6 = $ws.Range("A6:E6").Value2

$ws.Range("A5:E5").Value2 = $row6
$ws.Range("A6:E6").Value2 = $row5
